$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-03-31 Sunday" "2024-04-01 Monday"

Replace-Text "84×67=5628" "34×17=578"
Replace-Text "39×97=3783" "91×39=3549"
Replace-Text "17×95=1615" "40×68=2720"
Replace-Text "94×29=2726" "93×27=2511"
Replace-Text "42×71=2982" "47×19=893"

Replace-Text "88×21=1848" "18×24=432"
Replace-Text "47×16=752" "46×53=2438"
Replace-Text "31×77=2387" "68×16=1088"
Replace-Text "14×87=1218" "16×79=1264"
Replace-Text "65×11=715" "71×41=2911"

Replace-Text "90×91=8190" "61×31=1891"
Replace-Text "68×58=3944" "85×35=2975"
Replace-Text "83×92=7636" "70×23=1610"
Replace-Text "37×65=2405" "20×95=1900"
Replace-Text "57×39=2223" "82×82=6724"

Replace-Text "11×53=583" "73×34=2482"
Replace-Text "54×80=4320" "38×78=2964"
Replace-Text "22×14=308" "87×88=7656"
Replace-Text "11×36=396" "81×76=6156"
Replace-Text "56×38=2128" "41×47=1927"

Replace-Text "62×52=3224" "72×39=2808"
Replace-Text "29×86=2494" "60×48=2880"
Replace-Text "90×55=4950" "60×49=2940"
Replace-Text "46×19=874" "96×11=1056"
Replace-Text "92×50=4600" "66×25=1650"
